$d = $word.ActiveDocument

# Remove the stray leftover paragraph "?@tbl-chem-info" (an unresolved
# Quarto cross-reference placeholder) that was left behind in the
# "Other analyses" section. Delete the whole paragraph, including its
# trailing paragraph mark, so the surrounding text collapses back
# together exactly as it was before the placeholder was added.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*?@tbl-chem-info*") {
        $p.Range.Delete()
    }
}
